# Auto-generated Excel COM-interop edit script
# Applies the Tue Oct 24 21:41:17 UTC 2023 cryptos-list refresh:
# updates Price (D) / Volume(1h) (E) figures and restores the
# Aave / VeChain row order+data (rows 39-40) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold plain text that looks numeric (e.g. "1.00",
# "33.610.34"). Assigning .Value directly would let Excel coerce these to
# real numbers, so force text format first and strip it again afterwards
# so the cells end up with no explicit style, matching the source file.
# (Use one contiguous range, not a comma-union, so the format reliably
# reaches every row in the data table.)
$dPriceCells = "D2:D51"
$ws.Range($dPriceCells).NumberFormat = "@"

$ws.Range("D2").Value = "33.610.34"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "1.773.45"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "223.87"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "30.02"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "46.45"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").Value = "0.276"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "0.0659"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("D14").Value = "1.776.02"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "0.620"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "33.653.33"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").Value = "10.05"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "68.18"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "249.19"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "0.0₃0733"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "10.23"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "4.14"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "158.14"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "16.35"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "6.91"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "3.77"
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("D32").Value = "0.0512"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "3.53"
$ws.Range("E34").Value = "  +3.26%  "
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("D36").Value = "1.478.01"
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "0.625"
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0184"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "82.69"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "0.883"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "11.63"
$ws.Range("E50").Value = "  +12.42%  "
$ws.Range("D51").Value = "50.93"
$ws.Range("E51").Value = "  -3.17%  "

$ws.Range($dPriceCells).ClearFormats()

Write-Host "Updated 87 cell(s) in $($ws.Name) from the GitHub Actions cryptos refresh"
